$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.590.39"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.113.46"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").Value = "'1.011"
$ws.Range("E4").Value = "  +0.74%  "
$ws.Range("D5").Value = "'350.59"
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("D7").Value = "'0.5253"
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").Value = "'0.4515"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("D9").Value = "'54.44"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("D10").Value = "'0.09024"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").Value = "'1.174"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").Value = "'24.50"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "2.107.93"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").Value = "'6.819"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "'8.050"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "'101.49"
$ws.Range("E16").Value = "  +5.05%  "
$ws.Range("D17").Value = "'0.00001169"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").Value = "'1.012"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "'0.06724"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").Value = "'19.42"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").Value = "'6.292"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").Value = "30.647.57"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("E24").Value = "  +3.03%  "
$ws.Range("D25").Value = "'2.391"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("D26").Value = "2.361.12"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "'22.43"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "'164.98"
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("D29").Value = "'2.539"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").Value = "'136.79"
$ws.Range("E30").Value = "  +3.01%  "
$ws.Range("E31").Value = "  -4.31%  "
$ws.Range("D32").Value = "'0.1077"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").Value = "'1.664"
$ws.Range("E33").Value = "  -3.32%  "
$ws.Range("D34").Value = "'6.364"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("D35").Value = "'4.016"
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("D36").Value = "'10.38"
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("D37").Value = "'5.915"
$ws.Range("E37").Value = "  +5.99%  "
$ws.Range("E38").Value = "  +2.02%  "
$ws.Range("D39").Value = "'0.06843"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "'0.2311"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("D42").Value = "'0.6879"
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("D43").Value = "'1.271"
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("D44").Value = "'14.67"
$ws.Range("E44").Value = "  +4.14%  "
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").Value = "'0.6451"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D47").Value = "'3.759"
$ws.Range("E47").Value = "  +2.69%  "
$ws.Range("D48").Value = "'0.00000000361"
$ws.Range("E48").Value = "  +3.31%  "
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").Value = "'0.07290"
$ws.Range("D51").Value = "'82.38"
$ws.Range("E51").Value = "  -1.39%  "
